$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.170.23"
$ws.Range("E2").Value = "  +1.74%  "

# Row 3
$ws.Range("D3").Value = "2.778.47"
$ws.Range("E3").Value = "  +2.83%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'587.64"
$ws.Range("E5").Value = "  +0.52%  "

# Row 6
$ws.Range("D6").Value = "'161.89"
$ws.Range("E6").Value = "  +8.46%  "

# Row 7
$ws.Range("E7").Value = "  +2.14%  "

# Row 8
$ws.Range("E8").Value = "  +0.28%  "

# Row 9
$ws.Range("D9").Value = "2.787.77"
$ws.Range("E9").Value = "  +2.17%  "

# Row 10
$ws.Range("D10").Value = "'6.80"
$ws.Range("E10").Value = "  +1.19%  "

# Row 11
$ws.Range("D11").Value = "'0.115"
$ws.Range("E11").Value = "  +1.89%  "

# Row 12
$ws.Range("D12").Value = "'0.400"
$ws.Range("E12").Value = "  +3.46%  "

# Row 13
$ws.Range("E13").Value = "  +1.24%  "

# Row 14
$ws.Range("D14").Value = "3.271.43"
$ws.Range("E14").Value = "  +2.82%  "

# Row 15
$ws.Range("D15").Value = "'27.68"
$ws.Range("E15").Value = "  +3.80%  "

# Row 16
$ws.Range("D16").Value = "64.062.38"
$ws.Range("E16").Value = "  +1.78%  "

# Row 17
$ws.Range("E17").Value = "  +6.74%  "

# Row 18
$ws.Range("D18").Value = "2.780.97"
$ws.Range("E18").Value = "  +2.04%  "

# Row 19
$ws.Range("E19").Value = "  +4.28%  "

# Row 20
$ws.Range("D20").Value = "'5.06"
$ws.Range("E20").Value = "  +4.07%  "

# Row 21
$ws.Range("D21").Value = "'367.82"
$ws.Range("E21").Value = "  +1.62%  "

# Row 22
$ws.Range("D22").Value = "'7.08"
$ws.Range("E22").Value = "  +1.25%  "

# Row 23
$ws.Range("D23").Value = "'0.574"
$ws.Range("E23").Value = "  +8.47%  "

# Row 24
$ws.Range("E24").Value = "  +0.72%  "

# Row 25
$ws.Range("D25").Value = "'67.77"
$ws.Range("E25").Value = "  +3.59%  "

# Row 27
$ws.Range("D27").Value = "'8.88"
$ws.Range("E27").Value = "  +2.79%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0971"
$ws.Range("E28").Value = "  +13.65%  "

# Row 29
$ws.Range("E29").Value = "  +0.26%  "

# Row 30
$ws.Range("D30").Value = "'2.05"
$ws.Range("E30").Value = "  +1.42%  "

# Row 31
$ws.Range("D31").Value = "'7.31"
$ws.Range("E31").Value = "  +3.42%  "

# Row 32
$ws.Range("D32").Value = "'1.28"
$ws.Range("E32").Value = "  +6.74%  "

# Row 33
$ws.Range("D33").Value = "'172.56"
$ws.Range("E33").Value = "  +1.89%  "

# Row 34
$ws.Range("D34").Value = "'5.17"
$ws.Range("E34").Value = "  +8.83%  "

# Row 35
$ws.Range("D35").Value = "'20.88"
$ws.Range("E35").Value = "  +1.80%  "

# Row 36
$ws.Range("E36").Value = "  +0.12%  "

# Row 37
$ws.Range("D37").Value = "'1.51"
$ws.Range("E37").Value = "  +6.07%  "

# Row 38
$ws.Range("E38").Value = "  +2.43%  "

# Row 39
$ws.Range("E39").Value = "  +2.75%  "

# Row 40
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "'343.37"
$ws.Range("E40").Value = "  -1.94%  "

# Row 41
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'4.29"
$ws.Range("E41").Value = "  +0.90%  "

# Row 42
$ws.Range("D42").Value = "'6.31"
$ws.Range("E42").Value = "  +11.98%  "

# Row 43
$ws.Range("D43").Value = "'40.22"
$ws.Range("E43").Value = "  +2.59%  "

# Row 44
$ws.Range("D44").Value = "'22.53"
$ws.Range("E44").Value = "  +4.47%  "

# Row 45
$ws.Range("D45").Value = "'22.64"
$ws.Range("E45").Value = "  +4.67%  "

# Row 46
$ws.Range("D46").Value = "'0.0613"
$ws.Range("E46").Value = "  +3.54%  "

# Row 47
$ws.Range("D47").Value = "'0.653"
$ws.Range("E47").Value = "  +2.15%  "

# Row 48
$ws.Range("D48").Value = "'0.0263"
$ws.Range("E48").Value = "  +1.50%  "

# Row 49
$ws.Range("D49").Value = "'139.03"
$ws.Range("E49").Value = "  +0.38%  "

# Row 50
$ws.Range("E50").Value = "  +2.39%  "

# Row 51
$ws.Range("D51").Value = "2.181.26"
$ws.Range("E51").Value = "  +2.00%  "
